# Generate Report for Archive
#
# 1) Update status text "Ready for handoff" -> "In Translation" everywhere
#    it appears (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2) Shrink the corresponding "Status" columns to match the shorter text:
#    Overview columns E:F, zh-cn column C, de-de column C.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update the status values -------------------------------------------
if ($wsOverview.Range("E2").Text -eq $oldStatus) { $wsOverview.Range("E2").Value = $newStatus }
if ($wsOverview.Range("F2").Text -eq $oldStatus) { $wsOverview.Range("F2").Value = $newStatus }
if ($wsZhCn.Range("C2").Text -eq $oldStatus) { $wsZhCn.Range("C2").Value = $newStatus }
if ($wsDeDe.Range("C2").Text -eq $oldStatus) { $wsDeDe.Range("C2").Value = $newStatus }

# --- Resize the columns to fit the new, shorter text ---------------------
$newColumnWidth = 12.5

$wsOverview.Range("E1:F1").ColumnWidth = $newColumnWidth
$wsZhCn.Range("C1").ColumnWidth = $newColumnWidth
$wsDeDe.Range("C1").ColumnWidth = $newColumnWidth
